$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text (matches source data)
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2-8: refresh price (D) and volume (E) for unchanged coins
$ws.Range("D2").Value = '57.916.45'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '2.454.86'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '527.37'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '130.87'
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  +1.39%  '

# Rows 9-51: table shifted up one (LidoStakedEther dropped), refreshed coin data, Maker appended at bottom
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.0981'
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '4.96'
$ws.Range("E11").Value = '  -3.93%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '0.322'
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.890.89'
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '57.826.22'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '21.83'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.453.39'
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '10.36'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '4.15'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '316.69'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.07'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '64.98'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '0.408'
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '0.157'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '7.29'
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '172.17'
$ws.Range("E28").Value = '  +3.04%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0735'
$ws.Range("E29").Value = '  -2.60%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.70'
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '1.16'
$ws.Range("E31").Value = '  -3.09%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '6.13'
$ws.Range("E32").Value = '  -1.79%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '17.86'
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '1.18'
$ws.Range("E36").Value = '  -6.05%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '36.25'
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.46'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").Value = '0.805'
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.42'
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '267.93'
$ws.Range("E42").Value = '  -3.48%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '4.86'
$ws.Range("E43").Value = '  -4.50%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.585'
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '125.11'
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0931'
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = '0.0496'
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0212'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '17.06'
$ws.Range("E49").Value = '  -4.29%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '16.37'
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.712.06'
$ws.Range("E51").Value = '  -1.91%  '
